# "Implement add-images with excel data logic"
# The row for "Title 2" referenced an image (thumbnail2.jpg in column C /
# "thumbnail") that never got added, so its thumbnail reference is removed
# from the sheet (the row shrinks from 7 populated cells to 6 - column C is
# left blank for that row, while every other row keeps its thumbnail entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 == "Title 2" row. Clear the thumbnail cell (C3 = "thumbnail2.jpg").
$ws.Range("C3").Value2 = ""

# Reflect where the user was working when they made the change.
$ws.Range("C3").Select() | Out-Null
